$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = 5000
$ws.Range("F3").Value = 98

# Row 4
$ws.Range("C4").Value = 5000
$ws.Range("F4").Value = 98

# Row 5
$ws.Range("C5").Value = 5000
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 840
$ws.Range("F5").Value = 98
$ws.Range("G5").Value = 1

# Row 6
$ws.Range("C6").Value = 2000
$ws.Range("F6").Value = 85

# Row 7
$ws.Range("C7").Value = 1000
$ws.Range("F7").Value = 43

# Conditional formatting formula update on F3:F7 (97 -> 96 for both rules)
$fcs = $ws.Range("F3:F7").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $fc.Formula1 = "96"
}
